$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 10 (new training run "..._V2")
$ws.Range("A10").Value = "Inception-ResNet-V2_GRU_NoAttention_CustEmbedding_V2"
$ws.Range("B10").Value = "Inception-ResNet-V2 (512 units)"
$ws.Range("C10").Value = "GRU (512 units)"
$ws.Range("D10").Value = "20210503-163216"
$ws.Range("E10").Value = 0.439543646042204
$ws.Range("F10").Value = 0.535184995762384
$ws.Range("G10").Value = 0.678756956522537
$ws.Range("H10").Value = 0.857577816799581

# Clear out the cells in row 11 that are no longer populated
$ws.Range("E11").Clear()
$ws.Range("G11").Clear()

# Adjust column widths (COM ColumnWidth is quantized to whole pixels
# internally, so these are the closest settable values to the target
# stored widths of 26.28 / 15.7 / 16.58 characters)
$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(3).ColumnWidth = 14.8333333333333
$ws.Columns.Item(4).ColumnWidth = 15.6666666666667

# Move the active selection
$ws.Range("F15").Select()
